# Regenerate the randomized experiment task-order sheets.
# Tab positions (rId1..rId5) stay the same; only sheet names and the
# generated task-order rows change, per the new run's timestamped filenames.

$wb = $excel.ActiveWorkbook

# --- Rename tabs (positional; r:id -> worksheet mapping is unchanged) ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws5 = $wb.Worksheets.Item(5)

$ws1.Name = "vSAT_TO-16515890120514839"
$ws2.Name = "GNG_TO-1651589012082738"
$ws3.Name = "RS_TO-1651589012082738"
$ws4.Name = "NB_TO-1651589013968722"
$ws5.Name = "TOL_TO-16515890140155973"

# --- Sheet 1 (vSAT_TO...): same 4 data rows, new file names ---
$ws1.Range("B2").Value = "vSAT_stims-16515890120202353.csv"
$ws1.Range("B3").Value = "SAT_stims-1651589011988985.csv"
$ws1.Range("B4").Value = "vSAT_stims-165158901203586.csv"
$ws1.Range("B5").Value = "SAT_stims-16515890120046167.csv"

# --- Sheet 2 (GNG_TO...): shrinks from 10 rows to 5 rows ---
$ws2.Range("A6:A10").EntireRow.Delete()
$ws2.Range("B2").Value = "go_stims-16515890120514839.csv"
$ws2.Range("B3").Value = "GNG_stims-16515890120671113.csv"
$ws2.Range("B4").Value = "go_stims-16515890120671113.csv"
$ws2.Range("B5").Value = "GNG_stims-1651589012082738.csv"

# --- Sheet 3 (RS_TO...): unchanged content ---

# --- Sheet 4 (NB_TO...): grows from 7 rows to 10 rows ---
$ws4.Range("A7").Copy($ws4.Range("A8:A10"))
$ws4.Range("B2").Value = "OB-16515890127193604.csv"
$ws4.Range("B3").Value = "OB-16515890125766222.csv"
$ws4.Range("B4").Value = "ZB-match_2-16515890122022789.csv"
$ws4.Range("B5").Value = "ZB-match_5-16515890124025984.csv"
$ws4.Range("B6").Value = "TB-16515890130986774.csv"
$ws4.Range("B7").Value = "TB-16515890137624652.csv"
$ws4.Range("A8").Value = 6
$ws4.Range("B8").Value = "OB-1651589012924611.csv"
$ws4.Range("A9").Value = 7
$ws4.Range("B9").Value = "TB-16515890139531288.csv"
$ws4.Range("A10").Value = 8
$ws4.Range("B10").Value = "ZB-match_7-165158901217992.csv"

# --- Sheet 5 (TOL_TO...): grows from 5 rows to 7 rows ---
$ws5.Range("A5").Copy($ws5.Range("A6:A7"))
$ws5.Range("B2").Value = "MM_stims-16515890139843483.csv"
$ws5.Range("B3").Value = "ZM_stims-1651589013968722.csv"
$ws5.Range("B4").Value = "MM_stims-16515890139999723.csv"
$ws5.Range("B5").Value = "ZM_stims-16515890139843483.csv"
$ws5.Range("A6").Value = 4
$ws5.Range("B6").Value = "MM_stims-16515890140155973.csv"
$ws5.Range("A7").Value = 5
$ws5.Range("B7").Value = "ZM_stims-16515890139999723.csv"
